$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for rows 2, 3, 4 in the columns that change
$cols = @("A","B","E","F","G","H","Q","R")

$row2 = @{}
$row3 = @{}
$row4 = @{}

foreach ($c in $cols) {
    $row2[$c] = $ws.Range("${c}2").Value2
    $row3[$c] = $ws.Range("${c}3").Value2
    $row4[$c] = $ws.Range("${c}4").Value2
}

# Apply rotation: row2 -> row4, row3 -> row2, row4 -> row3
foreach ($c in $cols) {
    $ws.Range("${c}2").Value2 = $row3[$c]
    $ws.Range("${c}3").Value2 = $row4[$c]
    $ws.Range("${c}4").Value2 = $row2[$c]
}
